$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the species label on the Cebus row (A26): "Cebus imitator" -> "Cebus capucinus"
$ws.Range("A26").Value = "Cebus capucinus"

# Match the saved view/selection state (scroll position + active cell)
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F25").Select()
